# fix get dockercompose file
# Insert the missing repository rows (shift existing data down) and
# widen a couple of columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows, in final order (row 1 is the header row and stays untouched).
$data = @(
    @("alibghz/nestjs-microservices-docker", "Present", "Present", "Not present", "Not present", "Not present", "No CI/CD present", "Not present", "Not", "Not present"),
    @("sqshq/piggymetrics", "Present", "Present", "Not present", "Not present", "Present", "microservices touched by CI/CD", "Not present", "Not", "Present"),
    @("aidanwhiteley/books", "Present", "Present", "Not present", "Not present", "Not present", "No CI/CD present", "Not present", "Not", "Not present"),
    @("alanjeffares/notebook-to-microservice", "Present", "Present", "Not present", "Not present", "Not present", "No CI/CD present", "Not present", "Not", "Not present"),
    @("aliyun/alibabacloud-microservice-demo", "Present", "Present", "Not present", "Not present", "Not present", "microservices touched by CI/CD", "Not present", "Present", "Present")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# Widen columns A, F and J to accommodate the new data.
$ws.Columns.Item(1).ColumnWidth = 46
$ws.Columns.Item(6).ColumnWidth = 14.8
$ws.Columns.Item(10).ColumnWidth = 14.8
